$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Column C ("Förändrad") holds a date serial that changes from 45203 (2023-10-04)
# to 45204 (2023-10-05) for every data row (rows 2-97).
for ($r = 2; $r -le 97; $r++) {
    $cell = $ws.Cells.Item($r, 3)
    if ($cell.Value2 -eq 45203) {
        $cell.Value2 = 45204
    }
}
